$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows as of the final state: Vencimento (date serial), Taxa (flag), Data de Salvamento (timestamp)
$data = @(
    @(45792, 0, "2025-04-04 13:12:50"),
    @(45792, 1, "2025-04-04 13:13:02"),
    @(46249, 0, "2025-04-04 13:12:50"),
    @(46249, 1, "2025-04-04 13:13:02"),
    @(46522, 0, "2025-04-04 13:12:50"),
    @(46522, 1, "2025-04-04 13:13:02"),
    @(46980, 0, "2025-04-04 13:12:50"),
    @(46980, 1, "2025-04-04 13:13:02"),
    @(47253, 0, "2025-04-04 13:12:50"),
    @(47253, 1, "2025-04-04 13:13:02"),
    @(47710, 0, "2025-04-04 13:12:50"),
    @(47710, 1, "2025-04-04 13:13:02"),
    @(48441, 0, "2025-04-04 13:12:50"),
    @(48441, 1, "2025-04-04 13:13:02"),
    @(48714, 1, "2025-04-04 13:13:02"),
    @(48714, 0, "2025-04-04 13:12:50"),
    @(49444, 0, "2025-04-04 13:12:50"),
    @(49444, 1, "2025-04-04 13:13:02"),
    @(51363, 1, "2025-04-04 13:13:02"),
    @(51363, 0, "2025-04-04 13:12:50"),
    @(53097, 1, "2025-04-04 13:13:02"),
    @(53097, 0, "2025-04-04 13:12:50"),
    @(55015, 1, "2025-04-04 13:13:02"),
    @(55015, 0, "2025-04-04 13:12:50"),
    @(56749, 1, "2025-04-04 13:13:02"),
    @(56749, 0, "2025-04-04 13:12:50"),
    @(58668, 0, "2025-04-04 13:12:50"),
    @(58668, 1, "2025-04-04 13:13:02"),
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $rowData = $data[$i]
    $ws.Cells.Item($r, 1).Value2 = $rowData[0]
    $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 2).Value2 = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
}

